$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Column A: rename the existing test-user and add the new ones below it ---
$ws.Range("A2").Value = "DrGus"
$ws.Range("A3").Value = "GabrielaSaraiva"
$ws.Range("A4").Value = "GabrielaSantos"
$ws.Range("A5").Value = "MAriliaGabriela"
$ws.Range("A6").Value = "JoaoPedro"
$ws.Range("A7").Value = "Djonga"
$ws.Range("A8").Value = "Criolo"
$ws.Range("A9").Value = "BacudoExu"
$ws.Range("A10").Value = "SheldonNascimento"
$ws.Range("A11").Value = "ErikaBadu"

# --- Replicate the row-2 report formatting down through row 21 ---
for ($r = 3; $r -le 21; $r++) {
    $ws.Range("B$r").Style = "Hiperlink"
    $ws.Range("C$r").Style = "Hiperlink"
    $ws.Range("D$r").Style = "Hiperlink"
    $ws.Range("G$r").NumberFormat = "@"
    $ws.Range("K$r").NumberFormat = "@"
}

# row 12's A cell keeps the underline-only formatting that used to live on D4/G10
$ws.Range("A12").Font.Underline = 2

# --- Column A width grew slightly to fit the longer names ---
$ws.Columns.Item(1).ColumnWidth = 14.14

# --- Selection moved to A12 after the edits ---
$ws.Range("A12").Select()
